$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("A39").Value = 41221

$ws.Range("B39").Value = 3
$ws.Range("D39").Value = "Installer creation scripts continued, missing readMe files added"

$ws.Range("B39").Select()
